# Auto-generated edit script applying scheduled price-data refresh
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 31215.28
$ws.Range("I11").Value = 31215.28
$ws.Range("K11").Value = 31215.28
$ws.Range("M11").Value = -31075.28
$ws.Range("H55").Value = 137.35715
$ws.Range("I55").Value = 137.35715
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 137.35715
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 76.64285000000001
$ws.Range("N55").ClearContents()
$ws.Range("H64").Value = 3520.8718
$ws.Range("I64").Value = 2813.7827
$ws.Range("K64").Value = 2813.7827
$ws.Range("M64").Value = -2565.7827
$ws.Range("H67").Value = 3520.8718
$ws.Range("I67").Value = 2813.7827
$ws.Range("K67").Value = 2813.7827
$ws.Range("M67").Value = -1955.7827
$ws.Range("H70").Value = 8158.2
$ws.Range("J70").Value = 8731.444
$ws.Range("L70").Value = 26194.332
$ws.Range("N70").Value = -26734.332
$ws.Range("H73").Value = 8158.2
$ws.Range("J73").Value = 8731.444
$ws.Range("L73").Value = 26194.332
$ws.Range("N73").Value = -28066.332
$ws.Range("H112").Value = 2869.8262
$ws.Range("J112").Value = 3147.6
$ws.Range("L112").Value = 9442.799999999999
$ws.Range("N112").Value = -11658.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 133515.67
$ws.Range("I32").Value = 143600.4
$ws.Range("K32").Value = 143600.4
$ws.Range("M32").Value = -143313.4
$ws.Range("H45").Value = 2015
$ws.Range("I45").Value = 2096.6667
$ws.Range("K45").Value = 2096.6667
$ws.Range("M45").Value = -1719.6667
$ws.Range("H61").Value = 8246.700000000001
$ws.Range("I61").Value = 8718.556
$ws.Range("K61").Value = 8718.556
$ws.Range("M61").Value = -8506.556
$ws.Range("H74").Value = 16426
$ws.Range("I74").Value = 2497
$ws.Range("K74").Value = 2497
$ws.Range("M74").Value = -1623
$ws.Range("H77").Value = 16426
$ws.Range("I77").Value = 2497
$ws.Range("K77").Value = 12485
$ws.Range("M77").Value = -8117
$ws.Range("H125").Value = 68999.5
$ws.Range("J125").Value = 68999.5
$ws.Range("L125").Value = 68999.5
$ws.Range("N125").Value = -78839.5
$ws.Range("H132").Value = 759376.5600000001
$ws.Range("I132").Value = 894325.9399999999
$ws.Range("K132").Value = 2682977.82
$ws.Range("M132").Value = -2680447.82
$ws.Range("H136").Value = 8246.700000000001
$ws.Range("I136").Value = 8718.556
$ws.Range("K136").Value = 26155.668
$ws.Range("M136").Value = -23605.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 1001
$ws.Range("I31").Value = 1001
$ws.Range("K31").Value = 1001
$ws.Range("M31").Value = -749
$ws.Range("H94").Value = 1909.5385
$ws.Range("I94").Value = 1909.5385
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1909.5385
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1458.5385
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 12335.412
$ws.Range("I134").Value = 7113.467
$ws.Range("K134").Value = 21340.401
$ws.Range("M134").Value = -18805.401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 41
$ws.Range("I7").Value = 36.083332
$ws.Range("K7").Value = 36.083332
$ws.Range("M7").Value = 76.916668
$ws.Range("H23").Value = 167533.33
$ws.Range("I23").Value = 250050
$ws.Range("J23").Value = 2500
$ws.Range("K23").Value = 250050
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = -249810
$ws.Range("N23").Value = -2980
$ws.Range("H27").Value = 167533.33
$ws.Range("I27").Value = 250050
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 250050
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -249858
$ws.Range("N27").Value = -2884
$ws.Range("H31").Value = 3891.0833
$ws.Range("I31").Value = 4599.7144
$ws.Range("K31").Value = 4599.7144
$ws.Range("M31").Value = -4304.7144
$ws.Range("H34").Value = 3891.0833
$ws.Range("I34").Value = 4599.7144
$ws.Range("K34").Value = 4599.7144
$ws.Range("M34").Value = -4397.7144
$ws.Range("H93").Value = 11133.333
$ws.Range("I93").Value = 11133.333
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 11133.333
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -9261.333000000001
$ws.Range("N93").ClearContents()
$ws.Range("H131").Value = 37375
$ws.Range("J131").Value = 37375
$ws.Range("L131").Value = 37375
$ws.Range("N131").Value = -47455
$ws.Range("H132").Value = 2559
$ws.Range("I132").Value = 2559
$ws.Range("K132").Value = 7677
$ws.Range("M132").Value = -5147
$ws.Range("H134").Value = 4608.6665
$ws.Range("I134").Value = 3861.1428
$ws.Range("K134").Value = 11583.4284
$ws.Range("M134").Value = -9048.428400000001
$ws.Range("H141").Value = 214360.19
$ws.Range("J141").Value = 226317.53
$ws.Range("L141").Value = 226317.53
$ws.Range("N141").Value = -236677.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4962.25
$ws.Range("I32").Value = 4750
$ws.Range("K32").Value = 14250
$ws.Range("M32").Value = -13967
$ws.Range("H44").Value = 5614
$ws.Range("I44").Value = 499.2
$ws.Range("J44").Value = 8455.556
$ws.Range("K44").Value = 1497.6
$ws.Range("L44").Value = 25366.668
$ws.Range("M44").Value = -1099.6
$ws.Range("N44").Value = -26162.668
$ws.Range("H47").Value = 42.5
$ws.Range("J47").Value = 35
$ws.Range("L47").Value = 105
$ws.Range("N47").Value = -967
$ws.Range("H48").Value = 947.5
$ws.Range("J48").Value = 395
$ws.Range("L48").Value = 1185
$ws.Range("N48").Value = -1685
$ws.Range("H50").Value = 202064.56
$ws.Range("J50").Value = 627123.1
$ws.Range("L50").Value = 1881369.3
$ws.Range("N50").Value = -1882331.3
$ws.Range("H53").Value = 202064.56
$ws.Range("J53").Value = 627123.1
$ws.Range("L53").Value = 1881369.3
$ws.Range("N53").Value = -1882331.3
$ws.Range("H55").Value = 56004036
$ws.Range("I55").Value = 210000200
$ws.Range("J55").Value = 5436.364
$ws.Range("K55").Value = 630000600
$ws.Range("L55").Value = 16309.092
$ws.Range("M55").Value = -630000423
$ws.Range("N55").Value = -16663.092
$ws.Range("H87").Value = 15490.546
$ws.Range("J87").Value = 21999.834
$ws.Range("L87").Value = 65999.50199999999
$ws.Range("N87").Value = -68495.50199999999
$ws.Range("H90").Value = 15490.546
$ws.Range("J90").Value = 21999.834
$ws.Range("L90").Value = 197998.506
$ws.Range("N90").Value = -210478.506
$ws.Range("H109").Value = 6838.067
$ws.Range("I109").Value = 1945.1666
$ws.Range("J109").Value = 10100
$ws.Range("K109").Value = 5835.4998
$ws.Range("L109").Value = 30300
$ws.Range("M109").Value = -4795.4998
$ws.Range("N109").Value = -32380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2500
$ws.Range("I41").Value = 2500
$ws.Range("K41").Value = 2500
$ws.Range("M41").Value = -2145
$ws.Range("H126").Value = 9654.333000000001
$ws.Range("I126").Value = 10423.625
$ws.Range("K126").Value = 31270.875
$ws.Range("M126").Value = -28800.875
$ws.Range("H132").Value = 6910.1304
$ws.Range("I132").Value = 4993.026
$ws.Range("K132").Value = 14979.078
$ws.Range("M132").Value = -12449.078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1523
$ws.Range("I22").Value = 284.8889
$ws.Range("J22").Value = 2142.0557
$ws.Range("K22").Value = 284.8889
$ws.Range("L22").Value = 2142.0557
$ws.Range("M22").Value = 10.11110000000002
$ws.Range("N22").Value = -2732.0557
$ws.Range("H26").Value = 20250
$ws.Range("I26").Value = 17875
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 17875
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -17580
$ws.Range("N26").Value = -25590
$ws.Range("H27").Value = 1523
$ws.Range("I27").Value = 284.8889
$ws.Range("J27").Value = 2142.0557
$ws.Range("K27").Value = 284.8889
$ws.Range("L27").Value = 2142.0557
$ws.Range("M27").Value = -177.8889
$ws.Range("N27").Value = -2356.0557
$ws.Range("H30").Value = 2235.3333
$ws.Range("I30").Value = 2235.3333
$ws.Range("K30").Value = 2235.3333
$ws.Range("M30").Value = -2127.3333
$ws.Range("H31").Value = 2823
$ws.Range("I31").Value = 245.44444
$ws.Range("J31").Value = 10555.667
$ws.Range("K31").Value = 245.44444
$ws.Range("L31").Value = 10555.667
$ws.Range("M31").Value = 2.555560000000014
$ws.Range("N31").Value = -11051.667
$ws.Range("H46").Value = 2162.4167
$ws.Range("J46").Value = 2941.0588
$ws.Range("L46").Value = 2941.0588
$ws.Range("N46").Value = -3317.0588
$ws.Range("H75").Value = 40000
$ws.Range("I75").Value = 40000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 40000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -39064
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 40000
$ws.Range("I78").Value = 40000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 120000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -115320
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 14246.167
$ws.Range("I51").Value = 14246.167
$ws.Range("K51").Value = 14246.167
$ws.Range("M51").Value = -13736.167
$ws.Range("H126").Value = 2374.6667
$ws.Range("I126").Value = 2012
$ws.Range("K126").Value = 6036
$ws.Range("M126").Value = -3566

Write-Host "Applied scheduled market-data refresh: $($wb.Worksheets.Count) sheets touched."
